$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1248.0635
$ws.Range("I15").Value = 1248.0635
$ws.Range("K15").Value = 3744.1905
$ws.Range("M15").Value = -3575.1905

$ws.Range("H98").Value = 8431.799999999999
$ws.Range("I98").Value = 9555.823
$ws.Range("J98").Value = 2062.3333
$ws.Range("K98").Value = 9555.823
$ws.Range("L98").Value = 2062.3333
$ws.Range("M98").Value = -8057.823
$ws.Range("N98").Value = -5058.3333

$ws.Range("H116").Value = 4676.8887
$ws.Range("I116").Value = 3947.5557
$ws.Range("J116").Value = 6135.5557
$ws.Range("K116").Value = 3947.5557
$ws.Range("L116").Value = 6135.5557
$ws.Range("M116").Value = -505.5556999999999
$ws.Range("N116").Value = -13019.5557

$ws.Range("H122").Value = 8431.799999999999
$ws.Range("I122").Value = 9555.823
$ws.Range("J122").Value = 2062.3333
$ws.Range("K122").Value = 28667.469
$ws.Range("L122").Value = 6186.999899999999
$ws.Range("M122").Value = -26217.469
$ws.Range("N122").Value = -11086.9999

$ws.Range("H125").Value = 3699.5
$ws.Range("I125").Value = 4000
$ws.Range("J125").Value = 3499.1667
$ws.Range("K125").Value = 36000
$ws.Range("L125").Value = 31492.5003
$ws.Range("M125").Value = -33540
$ws.Range("N125").Value = -36412.5003

$ws.Range("H135").Value = 1622.9584
$ws.Range("I135").Value = 957.94116
$ws.Range("J135").Value = 3238
$ws.Range("K135").Value = 8621.470439999999
$ws.Range("L135").Value = 29142
$ws.Range("M135").Value = -6086.470439999999
$ws.Range("N135").Value = -34212

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 50456.76
$ws.Range("I32").Value = 54473.727
$ws.Range("K32").Value = 54473.727
$ws.Range("M32").Value = -54186.727

$ws.Range("H45").Value = 2078.375
$ws.Range("I45").Value = 2130.6
$ws.Range("J45").Value = 1991.3334
$ws.Range("K45").Value = 2130.6
$ws.Range("L45").Value = 1991.3334
$ws.Range("M45").Value = -1753.6
$ws.Range("N45").Value = -2745.3334

$ws.Range("H103").Value = 32333
$ws.Range("J103").Value = 32333
$ws.Range("L103").Value = 32333
$ws.Range("N103").Value = -34677

$ws.Range("H105").Value = 103663.2
$ws.Range("J105").Value = 103663.2
$ws.Range("L105").Value = 103663.2
$ws.Range("N105").Value = -110651.2

$ws.Range("H132").Value = 3851246
$ws.Range("I132").Value = 6670356.5
$ws.Range("J132").Value = 7004.0454
$ws.Range("K132").Value = 20011069.5
$ws.Range("L132").Value = 21012.1362
$ws.Range("M132").Value = -20008539.5
$ws.Range("N132").Value = -26072.1362

$ws.Range("H138").Value = 122500
$ws.Range("J138").Value = 122500
$ws.Range("L138").Value = 122500
$ws.Range("N138").Value = -132780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2309.0527
$ws.Range("I94").Value = 3028.5715
$ws.Range("K94").Value = 3028.5715
$ws.Range("M94").Value = -2577.5715

$ws.Range("H99").Value = 754.5
$ws.Range("I99").Value = 754.5
$ws.Range("K99").Value = 754.5
$ws.Range("M99").Value = 743.5

$ws.Range("H134").Value = 4704
$ws.Range("I134").Value = 3176.4849
$ws.Range("J134").Value = 8064.533
$ws.Range("K134").Value = 9529.4547
$ws.Range("L134").Value = 24193.599
$ws.Range("M134").Value = -6994.4547
$ws.Range("N134").Value = -29263.599

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5585.56
$ws.Range("I31").Value = 3214.8
$ws.Range("K31").Value = 3214.8
$ws.Range("M31").Value = -2919.8

$ws.Range("H34").Value = 5585.56
$ws.Range("I34").Value = 3214.8
$ws.Range("K34").Value = 3214.8
$ws.Range("M34").Value = -3012.8

$ws.Range("H58").Value = 5298.0435
$ws.Range("I58").Value = 3874.8
$ws.Range("J58").Value = 6392.846
$ws.Range("K58").Value = 3874.8
$ws.Range("L58").Value = 6392.846
$ws.Range("M58").Value = -3671.8
$ws.Range("N58").Value = -6798.846

$ws.Range("H136").Value = 5298.0435
$ws.Range("I136").Value = 3874.8
$ws.Range("J136").Value = 6392.846
$ws.Range("K136").Value = 11624.4
$ws.Range("L136").Value = 19178.538
$ws.Range("M136").Value = -9074.400000000001
$ws.Range("N136").Value = -24278.538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 444
$ws.Range("I6").Value = 44
$ws.Range("J6").Value = 710.6667
$ws.Range("K6").Value = 132
$ws.Range("L6").Value = 2132.0001
$ws.Range("M6").Value = -19
$ws.Range("N6").Value = -2358.0001

$ws.Range("H12").Value = 167.13637
$ws.Range("I12").Value = 8
$ws.Range("J12").Value = 202.5
$ws.Range("K12").Value = 24
$ws.Range("L12").Value = 607.5
$ws.Range("M12").Value = 149
$ws.Range("N12").Value = -953.5

$ws.Range("H40").Value = 2859.2856
$ws.Range("I40").Value = 35
$ws.Range("J40").Value = 3330
$ws.Range("K40").Value = 140
$ws.Range("L40").Value = 13320
$ws.Range("M40").Value = -71
$ws.Range("N40").Value = -13458

$ws.Range("H86").Value = 286.55554
$ws.Range("I86").Value = 396
$ws.Range("J86").Value = 199
$ws.Range("K86").Value = 1188
$ws.Range("L86").Value = 597
$ws.Range("M86").Value = -2
$ws.Range("N86").Value = -2969

$ws.Range("H89").Value = 286.55554
$ws.Range("I89").Value = 396
$ws.Range("J89").Value = 199
$ws.Range("K89").Value = 3564
$ws.Range("L89").Value = 1791
$ws.Range("M89").Value = 2364
$ws.Range("N89").Value = -13647

$ws.Range("H136").Value = 4691.5
$ws.Range("I136").Value = 2350
$ws.Range("K136").Value = 7050
$ws.Range("M136").Value = -1950

$ws.Range("H140").Value = 1951
$ws.Range("J140").Value = 2055.077
$ws.Range("L140").Value = 6165.231000000001
$ws.Range("N140").Value = -16525.231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4750
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 5000
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -4730
$ws.Range("N70").Value = -4540

$ws.Range("H73").Value = 4750
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 5000
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -4064
$ws.Range("N73").Value = -5872

$ws.Range("H102").Value = 2776.2
$ws.Range("I102").Value = 2776.2
$ws.Range("K102").Value = 2776.2
$ws.Range("M102").Value = -1154.2

$ws.Range("H132").Value = 5116.5
$ws.Range("I132").Value = 3185.5908
$ws.Range("K132").Value = 9556.7724
$ws.Range("M132").Value = -7026.7724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 701.8
$ws.Range("I16").Value = 619.82355
$ws.Range("K16").Value = 619.82355
$ws.Range("M16").Value = -449.82355

$ws.Range("H30").Value = 2046.1111
$ws.Range("I30").Value = 54.57143
$ws.Range("J30").Value = 9016.5
$ws.Range("K30").Value = 54.57143
$ws.Range("L30").Value = 9016.5
$ws.Range("M30").Value = 53.42857
$ws.Range("N30").Value = -9232.5

$ws.Range("H59").Value = 40000
$ws.Range("J59").Value = 40000
$ws.Range("L59").Value = 40000
$ws.Range("N59").Value = -41308

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 27000
$ws.Range("J101").Value = 27000
$ws.Range("L101").Value = 27000
$ws.Range("N101").Value = -33490

$ws.Range("H126").Value = 4869.7295
$ws.Range("I126").Value = 4726.533
$ws.Range("K126").Value = 14179.599
$ws.Range("M126").Value = -11709.599

$ws.Range("H136").Value = 3787.0454
$ws.Range("J136").Value = 6949.375
$ws.Range("L136").Value = 20848.125
$ws.Range("N136").Value = -25948.125
